$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price values that look numeric stay as text (matches source formatting)
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

# Apply the updated cell values
$ws.Range('D2').Value = '64.623.48'
$ws.Range('E2').Value = '  -2.94%  '
$ws.Range('D3').Value = '3.436.43'
$ws.Range('E3').Value = '  -2.60%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '581.81'
$ws.Range('E5').Value = '  -4.20%  '
$ws.Range('D6').Value = '134.63'
$ws.Range('E6').Value = '  -6.37%  '
$ws.Range('B7').Value = 'LidoStakedEther'
$ws.Range('C7').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D7').Value = '3.436.45'
$ws.Range('E7').Value = '  -2.50%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.484'
$ws.Range('E9').Value = '  -5.81%  '
$ws.Range('D10').Value = '0.121'
$ws.Range('E10').Value = '  -7.72%  '
$ws.Range('D11').Value = '6.99'
$ws.Range('E11').Value = '  -9.24%  '
$ws.Range('D12').Value = '0.377'
$ws.Range('E12').Value = '  -7.93%  '
$ws.Range('D13').Value = '4.016.63'
$ws.Range('E13').Value = '  -2.70%  '
$ws.Range('D14').Value = '0.0000179'
$ws.Range('E14').Value = '  -8.03%  '
$ws.Range('D15').Value = '3.437.35'
$ws.Range('E15').Value = '  -2.33%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value = '26.23'
$ws.Range('E16').Value = '  -8.70%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '0.115'
$ws.Range('E17').Value = '  -1.60%  '
$ws.Range('D18').Value = '64.487.28'
$ws.Range('E18').Value = '  -2.88%  '
$ws.Range('D19').Value = '9.66'
$ws.Range('E19').Value = '  -11.08%  '
$ws.Range('D20').Value = '5.69'
$ws.Range('E20').Value = '  -7.84%  '
$ws.Range('D21').Value = '13.60'
$ws.Range('E21').Value = '  -7.14%  '
$ws.Range('D22').Value = '380.85'
$ws.Range('E22').Value = '  -10.12%  '
$ws.Range('D23').Value = '0.544'
$ws.Range('E23').Value = '  -7.79%  '
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '72.00'
$ws.Range('E25').Value = '  -6.84%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = '5.73'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').Value = '3.569.64'
$ws.Range('E27').Value = '  -2.72%  '
$ws.Range('E28').Value = '  -7.57%  '
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').Value = '7.14'
$ws.Range('E30').Value = '  -9.68%  '
$ws.Range('D31').Value = '8.07'
$ws.Range('E31').Value = '  -9.67%  '
$ws.Range('D32').Value = '2.19'
$ws.Range('E32').Value = '  -11.34%  '
$ws.Range('D33').Value = '3.447.41'
$ws.Range('E33').Value = '  -2.50%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').Value = '22.83'
$ws.Range('E35').Value = '  -5.79%  '
$ws.Range('E36').Value = '  -9.84%  '
$ws.Range('D37').Value = '170.19'
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').Value = '6.74'
$ws.Range('E38').Value = '  -10.77%  '
$ws.Range('E39').Value = '  -12.43%  '
$ws.Range('E40').Value = '  -10.92%  '
$ws.Range('E41').Value = '  -10.98%  '
$ws.Range('D42').Value = '0.0762'
$ws.Range('E42').Value = '  -6.95%  '
$ws.Range('D43').Value = '0.801'
$ws.Range('E43').Value = '  -6.54%  '
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '41.88'
$ws.Range('E45').Value = '  -7.67%  '
$ws.Range('D46').Value = '4.28'
$ws.Range('E47').Value = '  -10.13%  '
$ws.Range('D48').Value = '22.62'
$ws.Range('E49').Value = '  -1.58%  '
$ws.Range('D50').Value = '6.53'
$ws.Range('E50').Value = '  -7.56%  '
$ws.Range('D51').Value = '2.197.74'
$ws.Range('E51').Value = '  -4.90%  '
